$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 2100
$ws.Range("I31").Value = 1750
$ws.Range("K31").Value = 5250
$ws.Range("M31").Value = -5020
$ws.Range("H109").Value = 30000
$ws.Range("J109").Value = 30000
$ws.Range("L109").Value = 30000
$ws.Range("N109").Value = -32774
$ws.Range("H132").Value = 1668.5536
$ws.Range("I132").Value = 1883.4634
$ws.Range("J132").Value = 1081.1333
$ws.Range("K132").Value = 5650.3902
$ws.Range("L132").Value = 3243.3999
$ws.Range("M132").Value = -3120.3902
$ws.Range("N132").Value = -8303.3999
$ws.Range("H135").Value = 954.1905
$ws.Range("I135").Value = 566.4865
$ws.Range("J135").Value = 3823.2
$ws.Range("K135").Value = 5098.3785
$ws.Range("L135").Value = 34408.8
$ws.Range("M135").Value = -2563.3785
$ws.Range("N135").Value = -39478.8
$ws.Range("H137").Value = 771.75555
$ws.Range("I137").Value = 688.875
$ws.Range("J137").Value = 866.4761999999999
$ws.Range("K137").Value = 2066.625
$ws.Range("L137").Value = 2599.4286
$ws.Range("M137").Value = 483.375
$ws.Range("N137").Value = -7699.428599999999
$ws.Range("H138").Value = 1950.99
$ws.Range("I138").Value = 827.78845
$ws.Range("J138").Value = 3167.7917
$ws.Range("K138").Value = 2483.36535
$ws.Range("L138").Value = 9503.375100000001
$ws.Range("M138").Value = 2656.63465
$ws.Range("N138").Value = -19783.3751
$ws.Range("H141").Value = 2523.0256
$ws.Range("I141").Value = 907.92
$ws.Range("J141").Value = 5407.143
$ws.Range("K141").Value = 2723.76
$ws.Range("L141").Value = 16221.429
$ws.Range("M141").Value = 2456.24
$ws.Range("N141").Value = -26581.429

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1769.39
$ws.Range("I32").Value = 1694.2963
$ws.Range("J32").Value = 2089.5264
$ws.Range("K32").Value = 1694.2963
$ws.Range("L32").Value = 2089.5264
$ws.Range("M32").Value = -1407.2963
$ws.Range("N32").Value = -2663.5264
$ws.Range("H61").Value = 1117.8928
$ws.Range("I61").Value = 975.4706
$ws.Range("J61").Value = 1338
$ws.Range("K61").Value = 975.4706
$ws.Range("L61").Value = 1338
$ws.Range("M61").Value = -763.4706
$ws.Range("N61").Value = -1762
$ws.Range("H74").Value = 822.93335
$ws.Range("I74").Value = 739.6539
$ws.Range("J74").Value = 1364.25
$ws.Range("K74").Value = 739.6539
$ws.Range("L74").Value = 1364.25
$ws.Range("M74").Value = 134.3461
$ws.Range("N74").Value = -3112.25
$ws.Range("H77").Value = 822.93335
$ws.Range("I77").Value = 739.6539
$ws.Range("J77").Value = 1364.25
$ws.Range("K77").Value = 3698.2695
$ws.Range("L77").Value = 6821.25
$ws.Range("M77").Value = 669.7304999999997
$ws.Range("N77").Value = -15557.25
$ws.Range("H132").Value = 1294.8846
$ws.Range("I132").Value = 1143.875
$ws.Range("J132").Value = 3107
$ws.Range("K132").Value = 3431.625
$ws.Range("L132").Value = 9321
$ws.Range("M132").Value = -901.625
$ws.Range("N132").Value = -14381
$ws.Range("H136").Value = 1117.8928
$ws.Range("I136").Value = 975.4706
$ws.Range("J136").Value = 1338
$ws.Range("K136").Value = 2926.4118
$ws.Range("L136").Value = 4014
$ws.Range("M136").Value = -376.4117999999999
$ws.Range("N136").Value = -9114

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 4425.909
$ws.Range("I8").Value = 1097.8572
$ws.Range("J8").Value = 10250
$ws.Range("K8").Value = 1097.8572
$ws.Range("L8").Value = 10250
$ws.Range("M8").Value = -957.8571999999999
$ws.Range("N8").Value = -10530
$ws.Range("H134").Value = 15465.353
$ws.Range("I134").Value = 1131.3508
$ws.Range("J134").Value = 73825.21000000001
$ws.Range("K134").Value = 3394.0524
$ws.Range("L134").Value = 221475.63
$ws.Range("M134").Value = -859.0523999999996
$ws.Range("N134").Value = -226545.63

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2503.0784
$ws.Range("I31").Value = 2251.2974
$ws.Range("J31").Value = 3168.5
$ws.Range("K31").Value = 2251.2974
$ws.Range("L31").Value = 3168.5
$ws.Range("M31").Value = -1956.2974
$ws.Range("N31").Value = -3758.5
$ws.Range("H34").Value = 2503.0784
$ws.Range("I34").Value = 2251.2974
$ws.Range("J34").Value = 3168.5
$ws.Range("K34").Value = 2251.2974
$ws.Range("L34").Value = 3168.5
$ws.Range("M34").Value = -2049.2974
$ws.Range("N34").Value = -3572.5
$ws.Range("H58").Value = 1337.3636
$ws.Range("I58").Value = 1412.1111
$ws.Range("J58").Value = 1001
$ws.Range("K58").Value = 1412.1111
$ws.Range("L58").Value = 1001
$ws.Range("M58").Value = -1209.1111
$ws.Range("N58").Value = -1407
$ws.Range("H132").Value = 1493.6323
$ws.Range("I132").Value = 992.4167
$ws.Range("J132").Value = 2057.5
$ws.Range("K132").Value = 2977.2501
$ws.Range("L132").Value = 6172.5
$ws.Range("M132").Value = -447.2501000000002
$ws.Range("N132").Value = -11232.5
$ws.Range("H134").Value = 1107.6487
$ws.Range("I134").Value = 1111.4117
$ws.Range("J134").Value = 1099.3043
$ws.Range("K134").Value = 3334.2351
$ws.Range("L134").Value = 3297.9129
$ws.Range("M134").Value = -799.2351000000003
$ws.Range("N134").Value = -8367.912899999999
$ws.Range("H136").Value = 1337.3636
$ws.Range("I136").Value = 1412.1111
$ws.Range("J136").Value = 1001
$ws.Range("K136").Value = 4236.3333
$ws.Range("L136").Value = 3003
$ws.Range("M136").Value = -1686.3333
$ws.Range("N136").Value = -8103

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 721.2632
$ws.Range("I5").Value = 729.1429000000001
$ws.Range("K5").Value = 2187.4287
$ws.Range("M5").Value = -2075.4287
$ws.Range("H45").Value = 1146.2222
$ws.Range("I45").Value = 825
$ws.Range("K45").Value = 2475
$ws.Range("M45").Value = -1943
$ws.Range("H74").Value = 5269.727
$ws.Range("I74").Value = 1980
$ws.Range("J74").Value = 5598.7
$ws.Range("K74").Value = 5940
$ws.Range("L74").Value = 16796.1
$ws.Range("M74").Value = -4879
$ws.Range("N74").Value = -18918.1
$ws.Range("H77").Value = 5269.727
$ws.Range("I77").Value = 1980
$ws.Range("J77").Value = 5598.7
$ws.Range("K77").Value = 17820
$ws.Range("L77").Value = 50388.3
$ws.Range("M77").Value = -12516
$ws.Range("N77").Value = -60996.3
$ws.Range("H81").Value = 2004.875
$ws.Range("I81").Value = 1504.3334
$ws.Range("J81").Value = 2305.2
$ws.Range("K81").Value = 4513.0002
$ws.Range("L81").Value = 6915.599999999999
$ws.Range("M81").Value = -3390.0002
$ws.Range("N81").Value = -9161.599999999999
$ws.Range("H84").Value = 2004.875
$ws.Range("I84").Value = 1504.3334
$ws.Range("J84").Value = 2305.2
$ws.Range("K84").Value = 13539.0006
$ws.Range("L84").Value = 20746.8
$ws.Range("M84").Value = -7923.000599999999
$ws.Range("N84").Value = -31978.8
$ws.Range("H130").Value = 3800
$ws.Range("J130").Value = 6750
$ws.Range("L130").Value = 20250
$ws.Range("N130").Value = -30290
$ws.Range("H135").Value = 721.2632
$ws.Range("I135").Value = 729.1429000000001
$ws.Range("K135").Value = 6562.2861
$ws.Range("M135").Value = -4027.2861

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 4758.6665
$ws.Range("I5").Value = 3052
$ws.Range("J5").Value = 5100
$ws.Range("K5").Value = 3052
$ws.Range("L5").Value = 5100
$ws.Range("M5").Value = -2940
$ws.Range("N5").Value = -5324
$ws.Range("H70").Value = 4286.25
$ws.Range("I70").Value = 3858
$ws.Range("K70").Value = 3858
$ws.Range("M70").Value = -3588
$ws.Range("H73").Value = 4286.25
$ws.Range("I73").Value = 3858
$ws.Range("K73").Value = 3858
$ws.Range("M73").Value = -2922
$ws.Range("H82").Value = 18633.334
$ws.Range("J82").Value = 18633.334
$ws.Range("L82").Value = 18633.334
$ws.Range("N82").Value = -19399.334
$ws.Range("H85").Value = 18633.334
$ws.Range("J85").Value = 18633.334
$ws.Range("L85").Value = 18633.334
$ws.Range("N85").Value = -21285.334
$ws.Range("H122").Value = 16228660
$ws.Range("I122").Value = 19956742
$ws.Range("J122").Value = 12500579
$ws.Range("K122").Value = 59870226
$ws.Range("L122").Value = 37501737
$ws.Range("M122").Value = -59867776
$ws.Range("N122").Value = -37506637
$ws.Range("H132").Value = 1885.4717
$ws.Range("I132").Value = 1741.4324
$ws.Range("J132").Value = 2218.5625
$ws.Range("K132").Value = 5224.2972
$ws.Range("L132").Value = 6655.6875
$ws.Range("M132").Value = -2694.2972
$ws.Range("N132").Value = -11715.6875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 1765.4524
$ws.Range("I132").Value = 1685.3513
$ws.Range("J132").Value = 2358.2
$ws.Range("K132").Value = 5056.0539
$ws.Range("L132").Value = 7074.599999999999
$ws.Range("M132").Value = -2526.0539
$ws.Range("N132").Value = -12134.6
$ws.Range("H136").Value = 1805.3208
$ws.Range("I136").Value = 1008.4889
$ws.Range("J136").Value = 6287.5
$ws.Range("K136").Value = 3025.4667
$ws.Range("L136").Value = 18862.5
$ws.Range("M136").Value = -475.4666999999999
$ws.Range("N136").Value = -23962.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 492.16666
$ws.Range("I113").Value = 240
$ws.Range("J113").Value = 672.2857
$ws.Range("K113").Value = 720
$ws.Range("L113").Value = 2016.8571
$ws.Range("M113").Value = 1450
$ws.Range("N113").Value = -6356.8571
$ws.Range("H122").Value = 1116.7273
$ws.Range("I122").Value = 960.8
$ws.Range("J122").Value = 1246.6666
$ws.Range("K122").Value = 2882.4
$ws.Range("L122").Value = 3739.9998
$ws.Range("M122").Value = -432.3999999999996
$ws.Range("N122").Value = -8639.9998
$ws.Range("H132").Value = 1306.3182
$ws.Range("I132").Value = 1301.7894
$ws.Range("J132").Value = 1335
$ws.Range("K132").Value = 3905.3682
$ws.Range("L132").Value = 4005
$ws.Range("M132").Value = -1375.3682
$ws.Range("N132").Value = -9065
$ws.Range("H136").Value = 1270.88
$ws.Range("I136").Value = 1493.1177
$ws.Range("J136").Value = 798.625
$ws.Range("K136").Value = 4479.3531
$ws.Range("L136").Value = 2395.875
$ws.Range("M136").Value = -1929.3531
